$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Search Page row: Execution Flag NO -> YES
$ws.Range("C3").Value = "YES"

# Clear out stray empty-but-styled cells / leftover values
$ws.Range("E3").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("I6").ClearContents()

# Move the active selection
$ws.Range("D18").Select()
